$wb = $excel.ActiveWorkbook

# --- ALC (sheet1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1405
$ws.Range("I28").Value = 485
$ws.Range("J28").Value = 2693
$ws.Range("K28").Value = 485
$ws.Range("L28").Value = 2693
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = -3663
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H137").Value = 2397.9
$ws.Range("I137").Value = 1497.1428
$ws.Range("K137").Value = 4491.428400000001
$ws.Range("M137").Value = -1941.428400000001

# --- ARM (sheet2) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 946.8570999999999
$ws.Range("I5").Value = 858.6667
$ws.Range("K5").Value = 858.6667
$ws.Range("M5").Value = -746.6667
$ws.Range("H39").Value = 6000
$ws.Range("I39").Value = 6000
$ws.Range("K39").Value = 6000
$ws.Range("M39").Value = -5480
$ws.Range("H61").Value = 3020.8
$ws.Range("I61").Value = 3020.8
$ws.Range("K61").Value = 3020.8
$ws.Range("M61").Value = -2808.8
$ws.Range("H74").Value = 2128.1428
$ws.Range("I74").Value = 2128.1428
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2128.1428
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1254.1428
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2128.1428
$ws.Range("I77").Value = 2128.1428
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 10640.714
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -6272.714
$ws.Range("N77").ClearContents()
$ws.Range("H92").Value = 96995.28999999999
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H97").Value = 30304836
$ws.Range("I97").Value = 47621310
$ws.Range("K97").Value = 47621310
$ws.Range("M97").Value = -47620814
$ws.Range("H102").Value = 70333336
$ws.Range("I102").Value = 5500009.5
$ws.Range("K102").Value = 5500009.5
$ws.Range("M102").Value = -5498387.5
$ws.Range("H132").Value = 3236.5
$ws.Range("I132").Value = 3236.5
$ws.Range("K132").Value = 9709.5
$ws.Range("M132").Value = -7179.5
$ws.Range("H136").Value = 3020.8
$ws.Range("I136").Value = 3020.8
$ws.Range("K136").Value = 9062.400000000001
$ws.Range("M136").Value = -6512.400000000001

# --- BSM (sheet3) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 946.8570999999999
$ws.Range("I4").Value = 858.6667
$ws.Range("K4").Value = 858.6667
$ws.Range("M4").Value = -743.6667
$ws.Range("H64").Value = 730
$ws.Range("I64").Value = 730
$ws.Range("K64").Value = 730
$ws.Range("M64").Value = -505
$ws.Range("H67").Value = 730
$ws.Range("I67").Value = 730
$ws.Range("K67").Value = 730
$ws.Range("M67").Value = 50
$ws.Range("H86").Value = 1986.909
$ws.Range("I86").Value = 1986.909
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1986.909
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -863.9090000000001
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1986.909
$ws.Range("I89").Value = 1986.909
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9934.545
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4318.545
$ws.Range("N89").ClearContents()

# --- CRP (sheet4) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 296.5625
$ws.Range("I7").Value = 283.13333
$ws.Range("K7").Value = 283.13333
$ws.Range("M7").Value = -170.13333
$ws.Range("H60").Value = 17750
$ws.Range("I60").Value = 14000
$ws.Range("J60").Value = 29000
$ws.Range("K60").Value = 14000
$ws.Range("L60").Value = 29000
$ws.Range("M60").Value = -13489
$ws.Range("N60").Value = -30022
$ws.Range("H134").Value = 2443
$ws.Range("I134").Value = 2426.7144
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 7280.1432
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -4745.1432
$ws.Range("N134").Value = -12570

# --- CUL (sheet5) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 58627.06
$ws.Range("I4").Value = 1441.9
$ws.Range("K4").Value = 4325.700000000001
$ws.Range("M4").Value = -4213.700000000001
$ws.Range("H87").Value = 800
$ws.Range("I87").Value = 800
$ws.Range("K87").Value = 2400
$ws.Range("M87").Value = -1152
$ws.Range("H90").Value = 800
$ws.Range("I90").Value = 800
$ws.Range("K90").Value = 7200
$ws.Range("M90").Value = -960
$ws.Range("H112").Value = 3600
$ws.Range("J112").Value = 3600
$ws.Range("L112").Value = 10800
$ws.Range("N112").Value = -13016
$ws.Range("H113").Value = 750.2
$ws.Range("I113").Value = 834.3333
$ws.Range("J113").Value = 624
$ws.Range("K113").Value = 2502.9999
$ws.Range("L113").Value = 1872
$ws.Range("M113").Value = -332.9998999999998
$ws.Range("N113").Value = -6212

# --- GSM (sheet6) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 6000
$ws.Range("K122").Value = 18000
$ws.Range("M122").Value = -15550
$ws.Range("H132").Value = 1457.7142
$ws.Range("I132").Value = 1457.7142
$ws.Range("K132").Value = 4373.142599999999
$ws.Range("M132").Value = -1843.142599999999

# --- LTW (sheet7) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 250
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -840
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 250
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 250
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -464
$ws.Range("H35").Value = 8708
$ws.Range("I35").Value = 1562.625
$ws.Range("K35").Value = 1562.625
$ws.Range("M35").Value = -1226.625
$ws.Range("H68").Value = 5975
$ws.Range("J68").Value = 6950
$ws.Range("L68").Value = 6950
$ws.Range("N68").Value = -8448
$ws.Range("H71").Value = 5975
$ws.Range("J71").Value = 6950
$ws.Range("L71").Value = 34750
$ws.Range("N71").Value = -42238
$ws.Range("H100").Value = 4497.5
$ws.Range("I100").Value = 4497.5
$ws.Range("K100").Value = 4497.5
$ws.Range("M100").Value = -3956.5

# --- WVR (sheet8) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6594.25
$ws.Range("I126").Value = 5670.2
$ws.Range("J126").Value = 8134.3335
$ws.Range("K126").Value = 17010.6
$ws.Range("L126").Value = 24403.0005
$ws.Range("M126").Value = -14540.6
$ws.Range("N126").Value = -29343.0005
